# Commit: "Changed final format from csv to tsv..." + corrected meteo sector names.
# Corrects two mislabeled "Secteur" values in the Meteo sheet:
#   "Petit Lotu" -> "Mezzanu"
#   "Mortella"   -> "A Torra di Murtella"
# and updates the last active selection to B18.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Meteo")

# Column B ("Secteur") rows 10-13 were "Petit Lotu" -> now "Mezzanu"
$ws.Range("B10:B13").Value2 = "Mezzanu"

# Column B ("Secteur") rows 14-17 were "Mortella" -> now "A Torra di Murtella"
$ws.Range("B14:B17").Value2 = "A Torra di Murtella"

# Move / record the active selection at B18 (was K6)
[void]$ws.Range("B18").Select()
